$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column V ("With bow"), matching the style of column U's header
$ws.Cells.Item(1, 21).Copy($ws.Cells.Item(1, 22))
$ws.Cells.Item(1, 22).Value = "With bow"

# Row 2 gets a FALSE flag in the new column
$ws.Cells.Item(2, 22).Value = $false

# Row 3 is duplicated from row 2's data (hull_id stays 2) plus a TRUE flag
$ws.Range("B2:U2").Copy($ws.Range("B3:U3"))
$ws.Cells.Item(3, 22).Value = $true

$excel.CutCopyMode = 0

# Update selection to reflect the next empty row
$ws.Range("A4").Select() | Out-Null
